$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 349; this shifts the existing rows 349:366
# down to 350:367 and carries formatting down from row 348 (so D349 keeps
# the date number format already used by the column).
$ws.Rows.Item(349).Insert()

# Populate the newly inserted row 349 with the new weekly record.
$ws.Range("A349").Value = 6
$ws.Range("B349").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C349").Value = "Metropolitana"
$ws.Range("D349").Value = 44585
$ws.Range("E349").Value = 13
$ws.Range("F349").Value = 100112043
$ws.Range("G349").Value = "Pepino ensalada"
$ws.Range("H349").Value = "Sin especificar"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 150
$ws.Range("K349").Value = 11000
$ws.Range("L349").Value = 12000
$ws.Range("M349").Value = 11533
$ws.Range("N349").Value = "$/caja 60 unidades"
$ws.Range("O349").Value = "Región de Arica y Parinacota"
$ws.Range("P349").Value = 192
$ws.Range("Q349").Value = 60
$ws.Range("R349").Value = "Hortaliza"
